# "add the detect system" — Scene.xlsx update
# Appends new detection-related task strings (torch/crystalball/lighthouse)
# to several scenes' "small-probability task" (G) and "quest" (F) columns.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Scene")

# Row 5 (昏暗密林): small-probability task gains a crystalball clause
$ws.Range("G5").Value = "forestfire;10|crystalball;35"

# Row 11 (黄金海岸): small-probability task gains a lighthouse clause
$ws.Range("G11").Value = "met;30|lighthouse;60"

# Row 12 (孤岛遗迹): small-probability task gains a lighthouse clause
$ws.Range("G12").Value = "met;30|goblinhome;40|lighthouse;70"

# Row 14 (玲珑峰): small-probability task gains a lighthouse clause
$ws.Range("G14").Value = "met;30|icedream;25|lighthouse;30"

# Row 16 (玲珑峰隧道): quest list gains a torch clause
$ws.Range("F16").Value = "mushroom;1|torch;1"

# Row 20 (落潮小径): small-probability task gains a lighthouse clause
$ws.Range("G20").Value = "met;30|lighthouse;60"

# Row 28 (村落入口): new crystalball quest in the previously-empty G column
$ws.Range("G28").Value = "crystalball;55"

# Row 29 (村中心): new crystalball quest in the previously-empty G column
$ws.Range("G29").Value = "crystalball;55"

# Restore the author's recorded selection after the edits
$ws.Range("G12").Select()
